# Commit: "add a Bland-Altman test"
#
# The "experiment" sheet (sheet2) contained two side-by-side tables:
#   A1 label "Helena" over table A2:D38 (testID 1..3)
#   F1 label "Yufeng" over table F2:I38 (testID 1..3)
# They get merged into a single stacked table in A:D (testID 1..6,
# with the former "Yufeng" block renumbered testID+3), the old
# per-person label row is dropped, and the "experiment" tab becomes
# the active sheet/tab (it was "stablize time" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "experiment"
$ws1 = $wb.Worksheets.Item(1)  # "stablize time"

# 1) Move the second ("Yufeng") block down below the first ("Helena")
#    block, renumbering its testID (column F / future column A) by +3
#    so the two tests don't collide (1..3 -> 4..6).
$ws.Range("F3:I38").Copy($ws.Range("A39"))
for ($r = 39; $r -le 74; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 3
}

# 2) Clear out the now-redundant "Yufeng" block (label, header, data).
$ws.Range("F1:I38").Clear()

# 3) Drop the leftover label row (old row 1 held "Helena"/"Yufeng"
#    headers) - this shifts the merged table up so it starts at row 2.
$ws.Rows.Item(1).Delete()

# 4) Re-create the column headers on the now-empty row 1.
$ws.Range("A1").Value = "observation"
$ws.Range("B1").Value = "time(s)"
$ws.Range("C1").Value = "watch"
$ws.Range("D1").Value = "device"

# 5) The "experiment" sheet becomes the active tab/sheet, with a new
#    selection; "stablize time" loses its selected/active tab mark
#    (its own selection is left untouched).
$ws.Activate()
$ws.Range("E13").Select()

# 6) Best-effort reproduction of the stray sheetFormatPr
#    outlineLevelCol="3" left on this sheet (column F already carries
#    explicit width formatting, so this only adds one attribute there
#    instead of fabricating a brand-new <col> entry).
$ws.Columns.Item(6).OutlineLevel = 3
